$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (col G) for rows 2-3
$wsOverview.Range("G2").Value = "2016-08-26 20:15:40"
$wsOverview.Range("G3").Value = "2016-08-26 20:15:40"

# zh-cn sheet: Priority (col E) ht -> mt, Correspond Handoff Datetime (col H),
# Correspond Handback DateTime (col K) for rows 2-3
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-26 20:15:35"
$wsZhCn.Range("H3").Value = "2016-08-26 20:15:35"
$wsZhCn.Range("K2").Value = "2016-08-26 20:15:52"
$wsZhCn.Range("K3").Value = "2016-08-26 20:15:52"

# de-de sheet: Priority (col E) ht -> mt, Correspond Handoff Datetime (col H),
# Correspond Handback DateTime (col K) for rows 2-3
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-26 20:15:40"
$wsDeDe.Range("H3").Value = "2016-08-26 20:15:40"
$wsDeDe.Range("K2").Value = "2016-08-26 20:16:01"
$wsDeDe.Range("K3").Value = "2016-08-26 20:16:01"
